# Applies the "minor changes" commit to the user manual:
#   1. Strike-through the "#5 - View Course Sequence (student)" bullet.
#   2. Move the hidden "_GoBack" bookmark from the very end of the document
#      to just after the bold "#9 - Drop Course (student)" heading
#      (this is what Word does automatically: it marks the last edited
#      spot with the _GoBack bookmark).
#   3. Justify ("both") a handful of body paragraphs that previously had
#      default (left) alignment.

$d = $word.ActiveDocument
$enDash = [char]0x2013
$wdAlignParagraphJustify = 3

# ---------------------------------------------------------------------
# 1. Strike-through "#5 - View Course Sequence (student)"
#    (use the Paragraphs collection directly so the paragraph mark's
#    own run properties pick up the formatting too, matching Word's
#    normal "select paragraph, press Ctrl+Shift+X" behaviour.)
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "#5*View Course Sequence (student)*") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------
# 2. Relocate the "_GoBack" bookmark
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the second (bold) "#9 - Drop Course (student)" heading.
$target = $d.Content
$target.Find.Execute("#9 $enDash Drop Course (student)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.Find.Execute("#9 $enDash Drop Course (student)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Insert a unique marker right after the heading text, then bookmark the
# marker and delete it again -- this avoids a COM quirk where adding a
# bookmark directly at a position that sits right on a paragraph boundary
# snaps to the wrong location.
$markerTag = [char]1 + "GOBACKMARK" + [char]1
$insertion = $d.Range($target.End, $target.End)
$insertion.InsertAfter($markerTag)

$markerRange = $d.Content
$markerRange.Find.Execute($markerTag, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Content
$markerRange2.Find.Execute($markerTag, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange2.Text = ""

# ---------------------------------------------------------------------
# 3. Justify several body paragraphs ("both" alignment)
# ---------------------------------------------------------------------
$justifyTargets = @(
    "An administrator that is logged in",
    "Clicking the manage courses option will display the",
    "The advanced search setting will display two search bars",
    "Once a course is searched and found,",
    "The view option will produce a page that displays the details of the",
    "The edit function is selected by clicking the pencil icon",
    "The final operation that an administrator can do to manage a course"
)

foreach ($needle in $justifyTargets) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            $p.Alignment = $wdAlignParagraphJustify
        }
    }
}
